# Generate Report for Handoff
# Replaces the stale handback-report data (two files that were handed back
# and in sync) with fresh handoff-report data (two files ready for handoff,
# one of which is a brand-new file id) across the Overview / zh-cn / de-de
# sheets, and drops the now-unused "Latest Target File" / "Latest Handback
# File" columns (F/G) from the per-locale sheets.

$wb = $excel.ActiveWorkbook

$oldMd1 = "a8a25305-d45e-4370-9a91-f732f20bbebb.md"
$oldMd2 = "cd7831bf-1df1-47ac-b351-225d4c84c32e.md"

$newMd1 = "8d8d2085-fe47-481b-b662-0dc702c42ed3.md"
$newMd2 = "ffff398fff67-c961-4876-b775-6cd7871573cb.md"

$newStatus = "Ready for handoff"
$newHandoffDate = "2016-03-22 15:11:09"

$newXlfZh = "8d8d2085-fe47-481b-b662-0dc702c42ed3.bc702ab10f0cf6487cf672ab2645395a00ab4626.zh-cn.xlf"
$newXlfDe = "8d8d2085-fe47-481b-b662-0dc702c42ed3.bc702ab10f0cf6487cf672ab2645395a00ab4626.de-de.xlf"

$newHandoffDatetime = "2016-03-22 15:11:02"
$newHandbackDatetime = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A2").Value2 = $newMd1
$wsOv.Range("B2").Value2 = $newStatus
$wsOv.Range("C2").Value2 = $newStatus
$wsOv.Range("D2").Value2 = $newHandoffDate

$wsOv.Range("A3").Value2 = $newMd2
$wsOv.Range("B3").Value2 = $newStatus
$wsOv.Range("C3").Value2 = $newStatus
$wsOv.Range("D3").Value2 = $newHandoffDate

foreach ($h in $wsOv.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq $wsOv.Range("A2").Address()) {
        $h.TextToDisplay = $newMd1
    } elseif ($addr -eq $wsOv.Range("A3").Address()) {
        $h.TextToDisplay = $newMd2
    }
}

# ---------------------------------------------------------------------
# Sheets "zh-cn" and "de-de": same 12-column layout, only the xlf suffix
# and handoff-datetime value differ between the two locales.
# ---------------------------------------------------------------------
$locales = @("zh-cn", "de-de")
foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale)

    if ($locale -eq "zh-cn") {
        $newXlf = $newXlfZh
    } else {
        $newXlf = $newXlfDe
    }

    # Remember the addresses of the columns we are about to drop (F/G) so we
    # can remove their hyperlinks before clearing the cells.
    $dropAddrs = @($ws.Range("F2").Address(), $ws.Range("G2").Address(), $ws.Range("F3").Address(), $ws.Range("G3").Address())

    $max = 50
    $i = 0
    $deleted = 0
    while ($deleted -lt $dropAddrs.Count -and $i -lt $max) {
        $found = $false
        foreach ($h in $ws.Hyperlinks) {
            $addr = $h.Range.Address()
            if ($dropAddrs -contains $addr) {
                $h.Delete()
                $deleted = $deleted + 1
                $found = $true
                break
            }
        }
        if (-not $found) { break }
        $i = $i + 1
    }

    # Update the remaining hyperlinks' display text in place (keeps their
    # existing target addresses / relationship ids, same as the row 2 /
    # row 3 file-name and target-file cells they sit on).
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq $ws.Range("A2").Address()) {
            $h.TextToDisplay = $newMd1
        } elseif ($addr -eq $ws.Range("D2").Address()) {
            $h.TextToDisplay = $newXlf
        } elseif ($addr -eq $ws.Range("A3").Address()) {
            $h.TextToDisplay = $newMd2
        } elseif ($addr -eq $ws.Range("D3").Address()) {
            $h.TextToDisplay = $newXlf
        }
    }

    # Now drop the F/G cells (Latest Target File / Latest Handback File)
    # for both rows entirely.
    $ws.Range("F2:G3").Clear()

    # Row 2
    $ws.Range("A2").Value2 = $newMd1
    $ws.Range("C2").Value2 = $newStatus
    $ws.Range("D2").Value2 = $newXlf
    $ws.Range("H2").Value2 = $newHandbackDatetime

    # Row 3
    $ws.Range("A3").Value2 = $newMd2
    $ws.Range("C3").Value2 = $newStatus
    $ws.Range("D3").Value2 = $newXlf
    $ws.Range("H3").Value2 = $newHandbackDatetime

    if ($locale -eq "zh-cn") {
        $ws.Range("E2").Value2 = $newHandoffDatetime
        $ws.Range("E3").Value2 = $newHandoffDatetime
    } else {
        $ws.Range("E2").Value2 = $newHandoffDate
        $ws.Range("E3").Value2 = $newHandoffDate
    }
}
